# Fixed issue #13 Permitir que en los ficheros de metadatos dos columnas se
# puedan relacionar para crear SKOS jerarquicos.
#
# This inserts a new row (row 2) of short "slug" style codes for each of the
# existing 5 metadata columns, pushing the previous rows 2-4 down to rows
# 3-5, and corrects the xsd type of the "Numero de miembros del hogar"
# measure (column A, now row 5) from xsd:double to xsd:string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 2, 3, 4 down to rows 3, 4, 5 to make room for the new
# row of short codes.
$ws.Rows.Item(2).Insert()

# New row 2: short slug-style codes matching each header in row 1.
$ws.Range("A2").Value = "n-hogares"
$ws.Range("B2").Value = "numero-de-miembros-del-hogar"
$ws.Range("C2").Value = "municipio-codigo"
$ws.Range("D2").Value = "n-medio-de-miembros"
$ws.Range("E2").Value = "municipio-nombre"

# Row 5 (previously row 4): fix the xsd type for column A from xsd:double to
# xsd:string.
$ws.Range("A5").Value = "xsd:string"
